$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230; this shifts the existing rows 230-323
# down to 231-324 and leaves a blank row 230 (inheriting formatting from
# the row above, e.g. the date style on column D).
$ws.Rows.Item(230).Insert()

# Populate the new row 230 with the new data record.
$ws.Cells.Item(230, 1).Value  = 4
$ws.Cells.Item(230, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(230, 3).Value  = "Los Lagos"
$ws.Cells.Item(230, 4).Value  = 44636
$ws.Cells.Item(230, 5).Value  = 10
$ws.Cells.Item(230, 6).Value  = 100114013
$ws.Cells.Item(230, 7).Value  = "Zanahoria"
$ws.Cells.Item(230, 8).Value  = "Sin especificar"
$ws.Cells.Item(230, 9).Value  = "Primera"
$ws.Cells.Item(230, 10).Value = 150
$ws.Cells.Item(230, 11).Value = 10000
$ws.Cells.Item(230, 12).Value = 10000
$ws.Cells.Item(230, 13).Value = 10000
$ws.Cells.Item(230, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(230, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(230, 16).Value = 500
$ws.Cells.Item(230, 17).Value = 20
$ws.Cells.Item(230, 18).Value = "Hortaliza"
